# The deck's design theme ("Integral") is swapped back to the default
# Office Theme color palette (this is what the underlying
# ppt/theme/theme1.xml - the slide master's theme - ends up containing
# after the edit; ppt/theme/theme2.xml, used only by the Notes Master,
# already carries the Office Theme palette and is left alone).
#
# PowerPoint exposes the 12 theme colours (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) through Slide.ThemeColorScheme.Colors(1..12).RGB, in
# that exact order, and writes land on the shared slide-master theme
# part. RGB values use the COLORREF convention: val = R + G*256 + B*65536.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

# Target palette = the built-in "Office Theme" colours.
$officeThemeRgb = @(
    0x000000,  # 1  dk1       000000
    0xFFFFFF,  # 2  lt1       FFFFFF
    0x6A5444,  # 3  dk2       44546A (stored BGR: B,G,R -> 6A,54,44)
    0xE6E6E7,  # 4  lt2       E7E6E6 (BGR -> E6,E6,E7)
    0xD59B5B,  # 5  accent1   5B9BD5 (BGR -> D5,9B,5B)
    0x317DED,  # 6  accent2   ED7D31 (BGR -> 31,7D,ED)
    0xA5A5A5,  # 7  accent3   A5A5A5
    0x00C0FF,  # 8  accent4   FFC000 (BGR -> 00,C0,FF)
    0xC47244,  # 9  accent5   4472C4 (BGR -> C4,72,44)
    0x47AD70,  # 10 accent6   70AD47 (BGR -> 47,AD,70)
    0xC16305,  # 11 hlink     0563C1 (BGR -> C1,63,05)
    0x724F95   # 12 folHlink  954F72 (BGR -> 72,4F,95)
)

for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeRgb[$i - 1]
}

Write-Output "Applied Office Theme colour scheme to the slide master theme."
